# Rename the "populationsCSV" Property to "populationsFolder" in the
# ProjectConfiguration sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Find the row whose "Property" column (A) currently holds "populationsCSV"
# and rename it to "populationsFolder", leaving the Value/Description columns
# untouched.
$found = $ws.Cells.Find("populationsCSV", [Type]::Missing, [Type]::Missing, 1)
if ($found -ne $null) {
    $found.Value = "populationsFolder"
} else {
    $ws.Range("A7").Value = "populationsFolder"
}
